# Auto-generated edit script applying numeric corrections to H/I/J/K/L/M/N
# columns (currentAveragePrice, currentAveragePriceNQ/HQ, LevePriceNQ/HQ,
# LeveProfitNQ/HQ) across the ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 32999
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 32999
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 32999
$ws.Range("N3").Value = -33227
$ws.Range("H7").Value = 6475.3335
$ws.Range("I7").Value = 4450
$ws.Range("J7").Value = 10526
$ws.Range("K7").Value = 4450
$ws.Range("L7").Value = 10526
$ws.Range("M7").Value = -4338
$ws.Range("N7").Value = -10750
$ws.Range("H14").Value = 6475.3335
$ws.Range("I14").Value = 4450
$ws.Range("J14").Value = 10526
$ws.Range("K14").Value = 4450
$ws.Range("L14").Value = 10526
$ws.Range("M14").Value = -4259
$ws.Range("N14").Value = -10908
$ws.Range("H39").Value = 365.55
$ws.Range("I39").Value = 75.916664
$ws.Range("J39").Value = 800
$ws.Range("K39").Value = 227.749992
$ws.Range("L39").Value = 2400
$ws.Range("M39").Value = 68.25000800000001
$ws.Range("H51").Value = 7249.75
$ws.Range("I51").Value = 8999.5
$ws.Range("J51").Value = 5500
$ws.Range("K51").Value = 8999.5
$ws.Range("L51").Value = 5500
$ws.Range("M51").Value = -8515.5
$ws.Range("H97").Value = 790.3333
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 790.3333
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 2370.9999
$ws.Range("N97").Value = -3362.9999
$ws.Range("H100").Value = 1902.7059
$ws.Range("I100").Value = 795.1
$ws.Range("J100").Value = 3485
$ws.Range("K100").Value = 795.1
$ws.Range("L100").Value = 3485
$ws.Range("M100").Value = -254.1
$ws.Range("N100").Value = -4567
$ws.Range("H102").Value = 32999
$ws.Range("I102").Value = 0
$ws.Range("J102").Value = 32999
$ws.Range("K102").Value = 0
$ws.Range("L102").Value = 32999
$ws.Range("N102").Value = -39489
$ws.Range("H110").Value = 80000
$ws.Range("I110").Value = 0
$ws.Range("J110").Value = 80000
$ws.Range("K110").Value = 0
$ws.Range("L110").Value = 80000
$ws.Range("N110").Value = -88180
$ws.Range("H117").Value = 75000
$ws.Range("I117").Value = 0
$ws.Range("J117").Value = 75000
$ws.Range("K117").Value = 0
$ws.Range("L117").Value = 75000
$ws.Range("N117").Value = -84178

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3512.4
$ws.Range("I45").Value = 2063.3333
$ws.Range("J45").Value = 4133.4287
$ws.Range("K45").Value = 2063.3333
$ws.Range("L45").Value = 4133.4287
$ws.Range("M45").Value = -1686.3333
$ws.Range("N45").Value = -4887.4287
$ws.Range("H110").Value = 3481.7778
$ws.Range("I110").Value = 2803.5
$ws.Range("J110").Value = 4838.3335
$ws.Range("K110").Value = 2803.5
$ws.Range("L110").Value = 4838.3335
$ws.Range("M110").Value = -758.5
$ws.Range("N110").Value = -8928.333500000001
$ws.Range("H133").Value = 100261
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 100261
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 100261
$ws.Range("N133").Value = -105321

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 4536.6665
$ws.Range("I86").Value = 2338.8333
$ws.Range("J86").Value = 6734.5
$ws.Range("K86").Value = 2338.8333
$ws.Range("L86").Value = 6734.5
$ws.Range("M86").Value = -1215.8333
$ws.Range("H89").Value = 4536.6665
$ws.Range("I89").Value = 2338.8333
$ws.Range("J89").Value = 6734.5
$ws.Range("K89").Value = 11694.1665
$ws.Range("L89").Value = 33672.5
$ws.Range("M89").Value = -6078.166499999999
$ws.Range("H94").Value = 322.4
$ws.Range("I94").Value = 269.33334
$ws.Range("J94").Value = 800
$ws.Range("K94").Value = 269.33334
$ws.Range("L94").Value = 800
$ws.Range("M94").Value = 181.66666
$ws.Range("H99").Value = 2075.4285
$ws.Range("I99").Value = 2086.1667
$ws.Range("J99").Value = 2011
$ws.Range("K99").Value = 2086.1667
$ws.Range("L99").Value = 2011
$ws.Range("M99").Value = -588.1667000000002
$ws.Range("H105").Value = 2010
$ws.Range("I105").Value = 2010
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 2010
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = -263
$ws.Range("N105").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H39").Value = 1000
$ws.Range("I39").Value = 1000
$ws.Range("J39").Value = 0
$ws.Range("K39").Value = 1000
$ws.Range("L39").Value = 0
$ws.Range("M39").Value = -609
$ws.Range("H42").Value = 0
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 0
$ws.Range("K42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("M42").ClearContents()
$ws.Range("H49").Value = 1000
$ws.Range("I49").Value = 1000
$ws.Range("J49").Value = 0
$ws.Range("K49").Value = 1000
$ws.Range("L49").Value = 0
$ws.Range("M49").Value = -818
$ws.Range("H131").Value = 38374.5
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 38374.5
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 38374.5
$ws.Range("N131").Value = -48454.5
$ws.Range("H141").Value = 68675.2
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 68675.2
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 68675.2
$ws.Range("N141").Value = -79035.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 285938.56
$ws.Range("I4").Value = 333445
$ws.Range("J4").Value = 900
$ws.Range("K4").Value = 1000335
$ws.Range("L4").Value = 2700
$ws.Range("M4").Value = -1000223
$ws.Range("N4").Value = -2924
$ws.Range("H128").Value = 783332.7
$ws.Range("I128").Value = 783332.7
$ws.Range("J128").Value = 0
$ws.Range("K128").Value = 2349998.1
$ws.Range("L128").Value = 0
$ws.Range("M128").Value = -2345018.1
$ws.Range("H132").Value = 1020.6
$ws.Range("I132").Value = 461.6
$ws.Range("J132").Value = 1579.6
$ws.Range("K132").Value = 4154.400000000001
$ws.Range("L132").Value = 14216.4
$ws.Range("M132").Value = -1624.400000000001
$ws.Range("H141").Value = 2185.8333
$ws.Range("I141").Value = 2185.8333
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 6557.499899999999
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = -1377.499899999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 18231562
$ws.Range("I7").Value = 14625500
$ws.Range("J7").Value = 24001260
$ws.Range("K7").Value = 14625500
$ws.Range("L7").Value = 24001260
$ws.Range("M7").Value = -14625388
$ws.Range("N7").Value = -24001484
$ws.Range("H8").Value = 18231562
$ws.Range("I8").Value = 14625500
$ws.Range("J8").Value = 24001260
$ws.Range("K8").Value = 14625500
$ws.Range("L8").Value = 24001260
$ws.Range("M8").Value = -14625361
$ws.Range("N8").Value = -24001538

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6726.143
$ws.Range("I7").Value = 4274.25
$ws.Range("J7").Value = 9995.333000000001
$ws.Range("K7").Value = 4274.25
$ws.Range("L7").Value = 9995.333000000001
$ws.Range("M7").Value = -4162.25
$ws.Range("H22").Value = 803
$ws.Range("I22").Value = 799.3333
$ws.Range("J22").Value = 811.25
$ws.Range("K22").Value = 799.3333
$ws.Range("L22").Value = 811.25
$ws.Range("M22").Value = -504.3333
$ws.Range("N22").Value = -1401.25
$ws.Range("H27").Value = 803
$ws.Range("I27").Value = 799.3333
$ws.Range("J27").Value = 811.25
$ws.Range("K27").Value = 799.3333
$ws.Range("L27").Value = 811.25
$ws.Range("M27").Value = -692.3333
$ws.Range("N27").Value = -1025.25
$ws.Range("H40").Value = 5981.846
$ws.Range("I40").Value = 6063.6665
$ws.Range("J40").Value = 5000
$ws.Range("K40").Value = 6063.6665
$ws.Range("L40").Value = 5000
$ws.Range("M40").Value = -5927.6665
$ws.Range("H61").Value = 4405.7334
$ws.Range("I61").Value = 2011
$ws.Range("J61").Value = 7142.5713
$ws.Range("K61").Value = 2011
$ws.Range("L61").Value = 7142.5713
$ws.Range("M61").Value = -1809
$ws.Range("H68").Value = 4365.8237
$ws.Range("I68").Value = 2781.9
$ws.Range("J68").Value = 6628.5713
$ws.Range("K68").Value = 2781.9
$ws.Range("L68").Value = 6628.5713
$ws.Range("M68").Value = -2032.9
$ws.Range("H71").Value = 4365.8237
$ws.Range("I71").Value = 2781.9
$ws.Range("J71").Value = 6628.5713
$ws.Range("K71").Value = 13909.5
$ws.Range("L71").Value = 33142.85649999999
$ws.Range("M71").Value = -10165.5
$ws.Range("H101").Value = 26500
$ws.Range("I101").Value = 0
$ws.Range("J101").Value = 26500
$ws.Range("K101").Value = 0
$ws.Range("L101").Value = 26500
$ws.Range("N101").Value = -32990
$ws.Range("H113").Value = 4405.7334
$ws.Range("I113").Value = 2011
$ws.Range("J113").Value = 7142.5713
$ws.Range("K113").Value = 2011
$ws.Range("L113").Value = 7142.5713
$ws.Range("M113").Value = 159
$ws.Range("H122").Value = 1500
$ws.Range("I122").Value = 1500
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 4500
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -2050
$ws.Range("H126").Value = 6726.143
$ws.Range("I126").Value = 4274.25
$ws.Range("J126").Value = 9995.333000000001
$ws.Range("K126").Value = 12822.75
$ws.Range("L126").Value = 29985.999
$ws.Range("M126").Value = -10352.75
$ws.Range("H132").Value = 1909.6
$ws.Range("I132").Value = 1819.2
$ws.Range("J132").Value = 2000
$ws.Range("K132").Value = 5457.6
$ws.Range("L132").Value = 6000
$ws.Range("M132").Value = -2927.6

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H95").Value = 47500
$ws.Range("I95").Value = 0
$ws.Range("J95").Value = 47500
$ws.Range("K95").Value = 0
$ws.Range("L95").Value = 47500
$ws.Range("N95").Value = -52992
$ws.Range("H96").Value = 1844.5555
$ws.Range("I96").Value = 1867.1666
$ws.Range("J96").Value = 1799.3334
$ws.Range("K96").Value = 1867.1666
$ws.Range("L96").Value = 1799.3334
$ws.Range("M96").Value = -494.1666
$ws.Range("H101").Value = 37800
$ws.Range("I101").Value = 0
$ws.Range("J101").Value = 37800
$ws.Range("K101").Value = 0
$ws.Range("L101").Value = 37800
$ws.Range("N101").Value = -44290
$ws.Range("H107").Value = 929.6
$ws.Range("I107").Value = 849.5
$ws.Range("J107").Value = 1250
$ws.Range("K107").Value = 2548.5
$ws.Range("L107").Value = 3750
$ws.Range("M107").Value = -628.5
$ws.Range("N107").Value = -7590
